# Atualizado por script em 11-11-2023 08:45
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: swap the "match detail" columns (F..V) between two rows, leaving
# the leading index/meta columns (A..E) untouched.
# ---------------------------------------------------------------------------
function Swap-MatchColumns($ws, $rowA, $rowB) {
    for ($col = 6; $col -le 22; $col++) {
        $cellA = $ws.Cells.Item($rowA, $col)
        $cellB = $ws.Cells.Item($rowB, $col)
        $tmp = $cellA.Value2
        $cellA.Value2 = $cellB.Value2
        $cellB.Value2 = $tmp
    }
}

# ---------------------------------------------------------------------------
# Rows 16/17: the two fixtures were swapped (same kickoff date, teams and
# odds exchanged).
# ---------------------------------------------------------------------------
Swap-MatchColumns $ws 16 17

# ---------------------------------------------------------------------------
# Rows 37/38/39: the three fixtures were cyclically rotated
#   new37 = old38, new38 = old39, new39 = old37
# ---------------------------------------------------------------------------
for ($col = 6; $col -le 22; $col++) {
    $save37 = $ws.Cells.Item(37, $col).Value2
    $ws.Cells.Item(37, $col).Value2 = $ws.Cells.Item(38, $col).Value2
    $ws.Cells.Item(38, $col).Value2 = $ws.Cells.Item(39, $col).Value2
    $ws.Cells.Item(39, $col).Value2 = $save37
}

# ---------------------------------------------------------------------------
# New row 82: appended fixture (Cape Town City vs Royal AM).
# Copy row 81's formatting down first so A82/E82 inherit the right styles.
# ---------------------------------------------------------------------------
$ws.Range("A81:V81").Copy($ws.Range("A82:V82"))

$ws.Cells.Item(82, 1).Value2 = 81
$ws.Cells.Item(82, 2).Value2 = "south-africa"
$ws.Cells.Item(82, 3).Value2 = "premier-league"
$ws.Cells.Item(82, 4).Value2 = "2023-2024"
$ws.Cells.Item(82, 5).Value2 = 45240.77083333334
$ws.Cells.Item(82, 6).Value2 = "Cape Town City"
$ws.Cells.Item(82, 7).Value2 = 2
$ws.Cells.Item(82, 8).Value2 = "Royal AM"
$ws.Cells.Item(82, 9).Value2 = 0
$ws.Cells.Item(82, 10).Value2 = 1.69
$ws.Cells.Item(82, 11).Value2 = "08/11/2023 16:42"
$ws.Cells.Item(82, 12).Value2 = 1.63
$ws.Cells.Item(82, 13).Value2 = "10/11/2023 18:21"
$ws.Cells.Item(82, 14).Value2 = 3.42
$ws.Cells.Item(82, 15).Value2 = "08/11/2023 16:42"
$ws.Cells.Item(82, 16).Value2 = 3.66
$ws.Cells.Item(82, 17).Value2 = "10/11/2023 18:21"
$ws.Cells.Item(82, 18).Value2 = 5.72
$ws.Cells.Item(82, 19).Value2 = "08/11/2023 16:42"
$ws.Cells.Item(82, 20).Value2 = 6.2
$ws.Cells.Item(82, 21).Value2 = "10/11/2023 18:21"
$ws.Cells.Item(82, 22).Value2 = "https://www.betexplorer.com/football/south-africa/premier-league/cape-town-city-royal-am/bDgEzuDR/"
